# "add new results to MAP" - append the augmentation sweep results to the
# "augmented_normal" sheet, and record the "select whole row 1" selection
# state that was left behind on the "normal" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "normal" sheet: the author selected the header row (row 1) before
#    saving - record that selection state. Do this first so the final
#    Activate() below (on augmented_normal) is what sticks as the
#    workbook's displayed/active tab.
# ---------------------------------------------------------------------
$wsNormal = $wb.Worksheets.Item("normal")
$wsNormal.Rows.Item(1).Select()

# ---------------------------------------------------------------------
# 2. "augmented_normal" sheet: fill in the new augmentation-ratio sweep
#    table (header row + 5 data rows).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("augmented_normal")
$ws.Activate()

$headers = @("# imgs", "augmentation", "mAP50 train", "mAP50-95 train", "mAP50 val", "mAP50-95 val")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$data = @(
    @(350, 0.1, 0.70035297086829296, 0.29831424425049102, 0.69996472803546395, 0.29838620068163502),
    @(350, 0.2, 0.78923872125068495, 0.31617320922684,    0.78929014193262603, 0.31612969102160698),
    @(350, 0.3, 0.83084999237384405, 0.42006941518102597, 0.82940956337066996, 0.419667219059496),
    @(350, 0.4, 0.87583359313041798, 0.42719661182555801, 0.87602208880721399, 0.42343376030268098),
    @(350, 0.5, 0.95690476190476104, 0.57743818197682495, 0.95690476190476104, 0.57743998049110601)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowVals[$c]
    }
}

# Header row is bold, like the other two sheets in the workbook.
$ws.Range("A1:F1").Font.Bold = $true

# Column widths, roughly matching the other sheets' autofit widths.
$ws.Range("B1").ColumnWidth = 15.333333333333334
$ws.Range("C1").ColumnWidth = 11.666666666666666
$ws.Range("D1").ColumnWidth = 15
$ws.Range("E1").ColumnWidth = 11.833333333333332
$ws.Range("F1").ColumnWidth = 12.833333333333332

# Leave the same cell selected as in the saved workbook.
$ws.Range("F8").Select()
